$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprints")

# Sprint Goal text for each sprint row
$ws.Range("F4").Value = "Initial planning of idea + deciding whether to pick fire alarm system or car security system"
$ws.Range("F5").Value = "Inital planning of project requirements and start writing project requirements"
$ws.Range("F6").Value = "Raspberry Pi board testing and finalising project requirments "
$ws.Range("F8").Value = "Writing of the main function,initalising components speicifcally on the temperature and smoke sensor, Getting output of buzzer,LCD and LED"
$ws.Range("F9").Value = "Working on manual switch and communication with SCDF and headless mode"
$ws.Range("F10").Value = "Testing of components and any access work not done from the previous sprint weeks"

# Updated Sprint 1 / Sprint 2 start dates (pushes the formula-derived dates forward).
# Use the raw date serial via Value2 so the existing custom date NumberFormat/style
# on the cell (s="13") is left untouched.
$ws.Range("G4").Value2 = 45775
$ws.Range("G8").Value2 = 45845

# Sprint state: first three sprints are now Closed
$ws.Range("I4").Value = "Closed"
$ws.Range("I5").Value = "Closed"
$ws.Range("I6").Value = "Closed"

$wb.Application.Calculate()
